$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear out the 10 "unhappy case" customer records (rows 16-25) - contents only,
# so rows below keep their original row numbers (no shift-up).
$ws.Range("A16:J25").ClearContents()

# Match the resulting selection left behind in the saved view state.
$ws.Range("A24:J25").Select()
